# --- Update the 2024 Jiangxi con-calendar workbook -------------------------
# The same edit applies to two worksheets with identical layout/content:
#   展览 (Exhibitions) and 全部类型 (All types).
# Changes:
#   1) A handful of 'want to go' / price corrections on existing rows 2-19.
#   2) Two brand-new 2024-04-20 events inserted as rows 20-21, pushing the
#      former rows 20-27 down to rows 22-29 (dimension grows from I27 to I29).
#   3) A few value corrections on the rows that shifted down.

$wb = $excel.ActiveWorkbook

$sheetNames = @('展览', '全部类型')
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- 1) Scalar corrections on existing rows (2-19) ---
    $ws.Cells.Item(2,6).Value = 1862
    $ws.Cells.Item(6,6).Value = 13
    $ws.Cells.Item(7,6).Value = 1547
    $ws.Cells.Item(8,6).Value = 20
    $ws.Cells.Item(9,6).Value = 603
    $ws.Cells.Item(10,7).Value = '不可售'
    $ws.Cells.Item(19,6).Value = 3618

    # --- 2) Insert two new rows at position 20; old rows 20-27 shift to 22-29 ---
    $ws.Rows.Item(20).Insert()
    $ws.Rows.Item(20).Insert()

    # Copy the index-column format down onto the two freshly inserted rows
    # so A20/A21 keep the same style as the rest of column A.
    $ws.Cells.Item(19,1).Copy()
    $ws.Cells.Item(20,1).PasteSpecial(-4122)
    $ws.Cells.Item(21,1).PasteSpecial(-4122)
    $excel.CutCopyMode = 0

    # --- Renumber column A (the 0-based index) for rows 20-29 ---
    $ws.Cells.Item(20,1).Value = 19
    $ws.Cells.Item(21,1).Value = 20
    $ws.Cells.Item(22,1).Value = 21
    $ws.Cells.Item(23,1).Value = 22
    $ws.Cells.Item(24,1).Value = 23
    $ws.Cells.Item(25,1).Value = 24
    $ws.Cells.Item(26,1).Value = 25
    $ws.Cells.Item(27,1).Value = 26
    $ws.Cells.Item(28,1).Value = 27
    $ws.Cells.Item(29,1).Value = 28

    # --- New row 20: 南昌·晨啼漫拥二次元随机舞蹈派对 (2024-04-20) ---
    $ws.Cells.Item(20,2).Value = '''2024-04-20'
    $ws.Cells.Item(20,3).Value = '南昌·晨啼漫拥二次元随机舞蹈派对-热爱欢聚(免费活动)'
    $ws.Cells.Item(20,4).Value = '莲塘镇澄湖东路1111号 玺悦城生活广场'
    $ws.Cells.Item(20,5).Value = '2024.04.20 15:00-04.20 19:00'
    $ws.Cells.Item(20,6).Value = 1
    $ws.Cells.Item(20,7).Value = 30.99
    $ws.Cells.Item(20,8).Value = 'https://show.bilibili.com/platform/detail.html?id=83272'
    $ws.Cells.Item(20,9).Value = '//i1.hdslb.com/bfs/openplatform/202403/wZUteBVO1710507652186.png'

    # --- New row 21: 抚州·四月之约动漫游戏聚会 (2024-04-20) ---
    $ws.Cells.Item(21,2).Value = '''2024-04-20'
    $ws.Cells.Item(21,3).Value = '抚州·四月之约动漫游戏聚会'
    $ws.Cells.Item(21,4).Value = '迎宾大道688号 抚州万达广场'
    $ws.Cells.Item(21,5).Value = '2024.04.20 10:00-04.20 17:00'
    $ws.Cells.Item(21,6).Value = 0
    $ws.Cells.Item(21,7).Value = 20
    $ws.Cells.Item(21,8).Value = 'https://show.bilibili.com/platform/detail.html?id=83316'
    $ws.Cells.Item(21,9).Value = '//i0.hdslb.com/bfs/openplatform/202403/2A7apu3o1711082007471.jpeg'

    # --- 3) Value corrections on the rows that shifted down (now 22-29) ---
    $ws.Cells.Item(23,6).Value = 324   # 吉安·COMIC LIFE次元假日04
    $ws.Cells.Item(24,6).Value = 486   # 景德镇·第十四届瓷都ACG动漫游戏博览会
    $ws.Cells.Item(25,6).Value = 188   # 江西·广电·Unlimited Project 动漫游戏博览会
    $ws.Cells.Item(25,7).Value = 19.9
    $ws.Cells.Item(26,6).Value = 339   # 江西·第二十二届九江ACJJ国际动漫展
    $ws.Cells.Item(28,6).Value = 1376  # 江西·ShiningStaR数字互娱嘉年华
    $ws.Cells.Item(29,6).Value = 140   # 南昌·代号鸢盛花行only
}

